$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Settings")

# The workbook's D8 cell holds a formula referencing the external
# FixedIncome.xla add-in (qlSerializationPath). Replace it with a plain
# literal value so the external reference is no longer needed, then
# sever/remove the external link itself.
$ws.Range("D8").Value = "'C:\Users\erik\junk\"

$wb.BreakLink("/WorkGroup/IMI_Workbooks/Production/QLXL_R01030x/framework/addin/FixedIncome.xla", 1)
